$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark off the "BEN BAO DAM / BEN NHAN BAO DAM"
#    signature line and onto the start of the "Muc dich su dung:" paragraph,
#    deleting the old run of text that used to read "Muc dich su dung: ".
# ---------------------------------------------------------------------------

# Locate + delete the run of text "Mục đích sử dụng: " (the whole run is
# removed, the following {{mdsd}} run is left untouched).
$target = $d.Content
$foundTarget = $target.Find.Execute("Mục đích sử dụng: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundTarget) {
    $target.Text = ""
    # Insert the (now relocated) "_GoBack" bookmark exactly where the
    # deleted run used to start.
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Strip the old "_GoBack" bookmark out of the signature paragraph while
#    leaving its text/formatting untouched.  Bookmark.Delete() is not wired
#    up in this host, so the paragraph is rewritten via Range.InsertXML with
#    the bookmark markers omitted (same paraId/rsids/run formatting).
# ---------------------------------------------------------------------------

$signature = $d.Content
$foundSig = $signature.Find.Execute("BÊN BẢO ĐẢM", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSig) {
    $paraStart = $signature.Start

    $tail = $d.Range($paraStart, $d.Content.End)
    $tail.Find.Execute("BÊN NHẬN BẢO ĐẢM", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $paraEnd = $tail.End

    $paraRange = $d.Range($paraStart, $paraEnd)

    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="579D23FA" w14:textId="6FB9C768" w:rsidR="00367D8B" w:rsidRDefault="00367D8B" w:rsidP="00511F45"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">              BÊN BẢO ĐẢM                                                      BÊN NHẬN BẢO ĐẢM</w:t></w:r></w:p>'
    $paraRange.InsertXML($xmlFrag) | Out-Null
}
